# issue #5: stock data from json to db
# Sheet "股票" (Stocks, sheet index 5) gains:
#   - a new "category" column inserted between property_category(H) and date(old I)
#     -> old I/J/K (date/legislator_name/legislator_id) slide right to J/K/L
#   - two new trailing columns: source_file (M) and index (N)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

$lastRow = 11

# --- Header row (row 1) ---------------------------------------------------
# Shift the old I1/J1/K1 header labels right by one (date, legislator_name,
# legislator_id), then relabel the now-vacated I1 as "category". Using
# Copy + PasteSpecial(xlPasteAll) (rather than reading/rewriting .Value)
# keeps the original cell type intact - e.g. it stops the literal text
# "2013-12-26" from being re-interpreted as a date serial number.
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4104)   # xlPasteAll

$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4104)

$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4104)

$ws.Range("I1").Value = "category"

# New trailing headers: source_file (M1), index (N1) - copy header formatting
$ws.Range("K1").Copy()
$ws.Range("M1").PasteSpecial(-4104)
$ws.Range("M1").Value = "source_file"

$ws.Range("K1").Copy()
$ws.Range("N1").PasteSpecial(-4104)
$ws.Range("N1").Value = "index"

# --- Data rows (rows 2..11) ------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    # Shift old I/J/K (date, legislator_name, legislator_id) right to J/K/L,
    # carrying the data-row formatting along.
    $ws.Range("K$r").Copy()
    $ws.Range("L$r").PasteSpecial(-4104)

    $ws.Range("J$r").Copy()
    $ws.Range("K$r").PasteSpecial(-4104)

    $ws.Range("I$r").Copy()
    $ws.Range("J$r").PasteSpecial(-4104)

    # New category value for every stock row.
    $ws.Range("I$r").Value = "normal"

    # New trailing columns: source_file + index.
    $ws.Range("K$r").Copy()
    $ws.Range("M$r").PasteSpecial(-4104)
    $ws.Range("M$r").Value = "tmpc7fb1"

    $ws.Range("K$r").Copy()
    $ws.Range("N$r").PasteSpecial(-4104)
    $ws.Range("N$r").Value = $ws.Range("A$r").Value()
}

Write-Output "stock sheet updated"
